# Update the "dSF" column (F) values for the rows that were recalculated
# ("repull data, push all data, mean calculation").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -1
    5  = 2
    6  = -2
    8  = -2
    14 = -9
    15 = -2
    17 = -5
    21 = -3
    22 = -4
    23 = -1
    28 = -4
    30 = 2
    33 = 6
    34 = -4
    36 = -2
    37 = 6
    39 = 3
    40 = -1
    41 = -2
    46 = 0
    49 = -1
    50 = 0
    59 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
